$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: fill in the missing AUTHOR (column B) for the existing
#     "The Lonely Dinosaur" story rows (2168-2180), by Henrique Figora ---
$ws.Range("B2168:B2180").Value = "Henrique Figora"

# --- Part 2: append the newly parsed story "Always Listen" by Brenda Wright
#     as rows 2181-2214 (LINE-NUMBER, AUTHOR, TITLE, SENTENCE) ---
$author = "Brenda Wright"
$title = "Always Listen"

$sentences = @(
    'Many moons ago there was a folk tale passed down from generation to generation. ',
    'The story is about the Cho-Cho Man or the Boogie Man. ',
    'It all began with a young boy named Dakota. ',
    'Dakota''s mother had asked him to gather some firewood for the evening fire. ',
    'The nights were getting cooler. ',
    'Instead of doing the chore right away, Dakota figured he had enough sunlight left to do two things. ',
    'He would gather the firewood after he played with his friends Morgan, Joey, and Bo. ',
    'Not long after, his mother began to wonder why Dakota never came with any firewood. ',
    'Night had fallen fast and it was too hard to see the firewood. ',
    'Dakota realized what had happened. ',
    'He told his friends about his predicament and they said that they would help their friend. ',
    'Soon after, the young boys were gathering wood in the dark. ',
    'Dakota''s mother was worried that something had happened. ',
    'She went looking for Dakota. ',
    'Not long after, she heard the boys. ',
    'She had mixed feelings. ',
    'She was happy to see Dakota but she was also mad that he hadn''t listened to his mother. ',
    'She wanted to teach Dakota and his friends a lesson so she sneaked into the bush very quietly. ',
    'Then she picked up a twig and threw it. ',
    'The twig hit the ground making noise. ',
    'It caught the boys'' attention. ',
    'She threw another broken twig in the opposite direction. ',
    'The boys jumped. ',
    'The mother had a hard time holding back her laughter and a muffled moan sounded through her fingers. ',
    'The boys were definitely scared. ',
    'But Dakota was the brave one. ',
    'He ventured forward. ',
    'The mother moaned again; she threw another twig. ',
    'The boys were transfixed. ',
    'Once more the mother moaned in a low and deep voice, "Listen to your mother." The boys dropped everything and ran. ',
    'The mother could not hold back her laughter and the boys, hearing this, scurried home even faster. ',
    'After the mother had picked up some firewood, she proceeded home, where she found Dakota in tears. ',
    'Dakota exclaimed, "Mother, Mother I''m sorry I didn''t listen to you." His mother held him with love and she knew that he was truly sorry. ',
    'From then on when Dakota felt like not listening, the mother would mention, "The boogie-man will scare you if you misbehave."'
)

$startRow = 2181
for ($i = 0; $i -lt $sentences.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Formula = "=ROW()-1"
    $ws.Cells.Item($r, 2).Value = $author
    $ws.Cells.Item($r, 3).Value = $title
    $ws.Cells.Item($r, 5).Value = $sentences[$i]
}

# --- Part 3: restore the view state (scroll position / selection) captured
#     in the saved workbook ---
$excel.ActiveWindow.ScrollRow = 2195
$ws.Range("B2211").Select()
